# Rename the AHB-Diff sheet's column headers so the "_old"/"_new" suffixes
# become the respective format-version suffixes ("_FV2210" / "_FV2304"),
# wrap the data range in a native Excel Table ("Table1"), and freeze the
# header row - matching the authors' "adapt column header formatting to
# respective input file names" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the 21 header cells in row 1 --------------------------------
$headers = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into a proper Excel table (Table1) ------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U" + $lastRow), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (View > Freeze Panes > Freeze Top Row) ------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
